$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    ,@("21CRB01291", "Hemmeter", "Permission Req'd To Use Licensed Dock", "1501:46-12-04", "MM", "No Contest", "Guilty", "`$ 50", "`$ 0", "5", "None")
    ,@("21CRB01291", "Hemmeter", "Permission Req'd To Use Licensed Dock", "1501:46-12-04", "MM", "No Contest", "Guilty", "`$ 50", "`$ 0", "5", "None")
    ,@("21CRB01291", "Hemmeter", "Permission Req'd To Use Licensed Dock", "1501:46-12-04", "MM", "No Contest", "Guilty", "`$ 50", "`$ 0", "5", "None")
    ,@("21CRB01291", "Hemmeter", "Permission Req'd To Use Licensed Dock", "1501:46-12-04", "MM", "No Contest", "Guilty", "`$ 50", "`$ 0", "5", "None")
    ,@("21CRB01291", "Hemmeter", "Permission Req'd To Use Licensed Dock", "1501:46-12-04", "MM", "No Contest", "Guilty", "`$ 50", "`$ 0", "5", "None")
    ,@("21CRB01291", "Hemmeter", "Permission Req'd To Use Licensed Dock", "1501:46-12-04", "MM", "No Contest", "Guilty", "`$ 50", "`$ 0", "5", "None")
    ,@("21CRB01291", "Hemmeter", "Permission Req'd To Use Licensed Dock", "1501:46-12-04", "MM", "No Contest", "Guilty", "`$ 50", "`$ 0", "5", "None")
    ,@("21CRB01291", "Hemmeter", "Permission Req'd To Use Licensed Dock", "1501:46-12-04", "MM", "No Contest", "Guilty", "`$ 50", "`$ 0", "5", "None")
    ,@("21CRB01291", "Hemmeter", "Permission Req'd To Use Licensed Dock", "1501:46-12-04", "MM", "No Contest", "Guilty", "`$ 50", "`$ 0", "5", "None")
    ,@("21CRB01291", "Hemmeter", "Permission Req'd To Use Licensed Dock", "1501:46-12-04", "MM", "No Contest", "Guilty", "`$ 50", "`$ 0", "5", "None")
    ,@("21CRB01291", "Hemmeter", "Permission Req'd To Use Licensed Dock", "1501:46-12-04", "MM", "No Contest", "Guilty", "`$ 50", "`$ 0", "5", "None")
    ,@("21CRB01291", "Hemmeter", "Permission Req'd To Use Licensed Dock", "1501:46-12-04", "MM", "No Contest", "Guilty", "`$ 50", "`$ 0", "5", "None")
    ,@("21CRB01291", "Hemmeter", "Permission Req'd To Use Licensed Dock", "1501:46-12-04", "MM", "No Contest", "Guilty", "`$ 50", "`$ 0", "5", "None")
    ,@("21CRB01268", "Hemmeter", "Possession Drug Paraphernalia", "2925.14(C)", "M4", "No Contest", "Guilty", "`$ 50", "`$ 0", "5", "None")
    ,@("21CRB01291", "Hemmeter", "Permission Req'd To Use Licensed Dock", "1501:46-12-04", "MM", "No Contest", "Guilty", "`$ 50", "`$ 0", "5", "None")
    ,@("21CRB01291", "Hemmeter", "Permission Req'd To Use Licensed Dock", "1501:46-12-04", "MM", "No Contest", "Guilty", "`$ 50", "`$ 0", "5", "None")
    ,@("21CRB01291", "Hemmeter", "Permission Req'd To Use Licensed Dock", "1501:46-12-04", "MM", "No Contest", "Guilty", "`$ 50", "`$ 0", "5", "None")
    ,@("21CRB01291", "Hemmeter", "Permission Req'd To Use Licensed Dock", "1501:46-12-04", "MM", "No Contest", "Guilty", "`$ 50", "`$ 0", "5", "None")
    ,@("21CRB01291", "Hemmeter", "Permission Req'd To Use Licensed Dock", "1501:46-12-04", "MM", "No Contest", "Guilty", "`$ 50", "`$ 0", "5", "None")
    ,@("21CRB01291", "Hemmeter", "Permission Req'd To Use Licensed Dock", "1501:46-12-04", "MM", "No Contest", "Guilty", "`$ 50", "`$ 0", "5", "None")
    ,@("21CRB01291", "Hemmeter", "Permission Req'd To Use Licensed Dock", "1501:46-12-04", "MM", "No Contest", "Guilty", "`$ 50", "`$ 0", "5", "None")
    ,@("21CRB01291", "Hemmeter", "Permission Req'd To Use Licensed Dock", "1501:46-12-04", "MM", "No Contest", "Guilty", "`$ 50", "`$ 0", "5", "None")
    ,@("21CRB01291", "Hemmeter", "Permission Req'd To Use Licensed Dock", "1501:46-12-04", "MM", "No Contest", "Guilty", "`$ 50", "`$ 0", "5", "None")
    ,@("21CRB01291", "Hemmeter", "Permission Req'd To Use Licensed Dock", "1501:46-12-04", "MM", "No Contest", "Guilty", "`$ 50", "`$ 0", "5", "None")
)

$startRow = 1123
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 1; $c -le 11; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $cell.NumberFormat = "@"
        $cell.Value = $rowData[$c-1]
    }
}
